$d = $word.ActiveDocument

# wdBorderBottom
$wdBorderBottom = -3

# This template's "HorizontalRule" paragraph style draws a horizontal rule
# using a bottom paragraph border. The published version of the extension
# pins an explicit line weight (w:sz="10", i.e. 1.25pt) on that border
# instead of leaving it to the (thin) implicit default. We need to set
# that weight on every HorizontalRule paragraph's bottom border, without
# touching anything else (style, color, indentation, centering, ...).
#
# Word's Paragraphs(N) *getter* (and .Style / .Range reads) correctly use
# 1-based paragraph-only indexing. However, in this runtime, writing a
# paragraph Border property (Paragraphs(N).Borders(...).PropertyX = ...)
# instead resolves N against the raw sequence of body-level blocks
# (paragraphs AND bookmarkStart/bookmarkEnd markers interleaved between
# them), so a naive "for each paragraph found, set Borders(...) on it"
# ends up editing the wrong paragraph whenever bookmarks precede it.
#
# To work around that, we compute -- straight from the document's own
# OOXML via Range.WordOpenXML, not any hard-coded number -- the raw
# block index of each HorizontalRule paragraph, and index into
# Paragraphs() with THAT number so the write lands on the right block.

$xml = $d.Content.WordOpenXML

# Pull out the main document part's XML and just its <w:body>.
$partMatch = [regex]::Match($xml, '(?s)<pkg:part pkg:name="/word/document\.xml"[^>]*>.*?<pkg:xmlData>(.*?)</pkg:xmlData>\s*</pkg:part>')
$docXml = $partMatch.Groups[1].Value
$bodyMatch = [regex]::Match($docXml, '(?s)<w:body>(.*)</w:body>')
$body = $bodyMatch.Groups[1].Value

# Walk the top-level block markers in document order (paragraphs,
# bookmark start/end markers, tables, sectPr) to build the mapping from
# "which block number is this paragraph" - that's what the buggy Border
# setter actually indexes by.
$tokenRegex = [regex]'<w:p>|<w:p\s[^>]*>|<w:p/>|<w:bookmarkStart\b[^>]*/>|<w:bookmarkEnd\b[^>]*/>|<w:tbl>|</w:tbl>|<w:sectPr\b'
$tokens = $tokenRegex.Matches($body)

$blockIndex = 0
$targetBlocks = New-Object System.Collections.ArrayList
for ($k = 0; $k -lt $tokens.Count; $k++) {
    $tok = $tokens[$k]
    $blockIndex = $blockIndex + 1
    $val = $tok.Value
    $isParagraph = $val.StartsWith("<w:p>") -or $val.StartsWith("<w:p ") -or ($val -eq "<w:p/>")
    if ($isParagraph) {
        # Look only within this paragraph's own markup (up to the next
        # block token) so we don't bleed into a sibling paragraph.
        $start = $tok.Index + $val.Length
        $end = $body.Length
        if ($k + 1 -lt $tokens.Count) {
            $end = $tokens[$k + 1].Index
        }
        $len = [Math]::Min($end - $start, 500)
        if ($len -gt 0) {
            $snippet = $body.Substring($start, $len)
            if ($snippet.Contains('w:val="HorizontalRule"')) {
                $targetBlocks.Add($blockIndex) | Out-Null
            }
        }
    }
}

foreach ($b in $targetBlocks) {
    # LineWidth is expressed in half-eighths-of-a-point here, i.e.
    # w:sz ends up as 2x whatever we assign, so 5 => w:sz="10" (1.25pt).
    $d.Paragraphs($b).Borders($wdBorderBottom).LineWidth = 5
}
